# update scripts with new tpm
# Recompute the Wnt3-Fzd7 LR-pairs sheet: the old "ECs" sending-cluster rows
# are dropped, the former "FAPs" sending-cluster rows (recalculated with the
# new TPM values) move up to become rows 2-4, and the now-unused trailing
# rows are removed so the sheet shrinks from A1:T7 to A1:T4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: FAPs -> Wnt3 -> Fzd7 -> ECs -----------------------------------
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt3"
$ws.Range("C2").Value = "Fzd7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1741663333333333
$ws.Range("H2").Value = 0.5224989999999999
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 1.01111
$ws.Range("N2").Value = 3.03333
$ws.Range("O2").Value = 0.04063212692754557
$ws.Range("P2").Value = 0.04063212692754556
$ws.Range("Q2").Value = 0.1761013212966666
$ws.Range("R2").Value = 1.58491189167
$ws.Range("S2").Value = 0.04063212692754557
$ws.Range("T2").Value = 0.04063212692754556

# --- Row 3: FAPs -> Wnt3 -> Fzd7 -> FAPs -----------------------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt3"
$ws.Range("C3").Value = "Fzd7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1741663333333333
$ws.Range("H3").Value = 0.5224989999999999
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.4065982422683317
$ws.Range("P3").Value = 0.4065982422683317
$ws.Range("Q3").Value = 1.762213625391444
$ws.Range("R3").Value = 15.859922628523
$ws.Range("S3").Value = 0.4065982422683317
$ws.Range("T3").Value = 0.4065982422683317

# --- Row 4: FAPs -> Wnt3 -> Fzd7 -> MuSCs ----------------------------------
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt3"
$ws.Range("C4").Value = "Fzd7"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1741663333333333
$ws.Range("H4").Value = 0.5224989999999999
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("O4").Value = 0.5527696308041227
$ws.Range("P4").Value = 0.5527696308041226
$ws.Range("Q4").Value = 2.395726478479889
$ws.Range("R4").Value = 21.561538306319
$ws.Range("S4").Value = 0.5527696308041227
$ws.Range("T4").Value = 0.5527696308041226

# --- Drop the old rows 5-7 (their recalculated data now lives in 2-4) -----
$ws.Rows("5:7").Delete()
